$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.128.92"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").Value = "1.636.61"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.90"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.518"
$ws.Range("E6").Value = "  +2.27%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.254"
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.93"
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "1.864.98"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").Value = "1.632.51"
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.13"
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("E15").Value = "  +2.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.64"
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("D17").Value = "27.125.66"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("E18").Value = "  +1.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.91"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.84"
$ws.Range("E21").Value = "  +1.39%  "
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("E23").Value = "  +2.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.11"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.62"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.39"
$ws.Range("E27").Value = "  +2.00%  "
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.38"
$ws.Range("E32").Value = "  +1.49%  "
$ws.Range("E33").Value = "  +0.41%  "
$ws.Range("D34").Value = "1.297.94"
$ws.Range("E34").Value = "  +2.73%  "
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("E36").Value = "  +1.11%  "
$ws.Range("E37").Value = "  -0.40%  "
$ws.Range("E38").Value = "  +1.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.543"
$ws.Range("E39").Value = "  +1.94%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.806"
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("E42").Value = "  +6.30%  "
$ws.Range("E43").Value = "  -1.43%  "
$ws.Range("D44").Value = "1.775.88"
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.65"
$ws.Range("E45").Value = "  -0.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.34"
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("E47").Value = "  +0.57%  "
$ws.Range("D48").Value = "0.0₆0106"
$ws.Range("E48").Value = "  +0.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0512"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.65"
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0957"
$ws.Range("E51").Value = "  -0.32%  "
